$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-CellText 'D2' '59.917.23'
Set-CellText 'E2' '  +2.94%  '
Set-CellText 'D3' '2.422.61'
Set-CellText 'E3' '  +3.01%  '
Set-CellText 'E4' '  +0.01%  '
Set-CellText 'D5' '554.78'
Set-CellText 'E5' '  +2.75%  '
Set-CellText 'D6' '138.33'
Set-CellText 'E6' '  +1.77%  '
Set-CellText 'D7' '1.00'
Set-CellText 'E7' '  -0.04%  '
Set-CellText 'D8' '0.571'
Set-CellText 'E8' '  +1.07%  '
Set-CellText 'E9' '  +5.06%  '
Set-CellText 'D10' '5.82'
Set-CellText 'E10' '  +4.65%  '
Set-CellText 'D11' '0.359'
Set-CellText 'E11' '  +1.68%  '
Set-CellText 'E12' '  -2.11%  '
Set-CellText 'D13' '24.74'
Set-CellText 'E13' '  +3.90%  '
Set-CellText 'D14' '2.852.32'
Set-CellText 'E14' '  +2.96%  '
Set-CellText 'D15' '59.809.71'
Set-CellText 'E15' '  +2.82%  '
Set-CellText 'E16' '  +4.39%  '
Set-CellText 'D17' '2.466.66'
Set-CellText 'E17' '  +4.87%  '
Set-CellText 'E18' '  +6.86%  '
Set-CellText 'E19' '  +3.41%  '
Set-CellText 'D20' '334.26'
Set-CellText 'E20' '  +0.70%  '
Set-CellText 'D21' '6.90'
Set-CellText 'E21' '  +1.20%  '
Set-CellText 'D22' '0.999'
Set-CellText 'E22' '  -0.09%  '
Set-CellText 'D23' '64.54'
Set-CellText 'E23' '  +2.78%  '
Set-CellText 'E24' '  +1.04%  '
Set-CellText 'D25' '8.55'
Set-CellText 'E25' '  +0.65%  '
Set-CellText 'D26' '1.00'
Set-CellText 'E26' '  -0.17%  '
Set-CellText 'E27' '  -0.44%  '
Set-CellText 'E28' '  +7.06%  '
Set-CellText 'E29' '  +3.38%  '
Set-CellText 'D30' '170.86'
Set-CellText 'E30' '  -0.22%  '
Set-CellText 'D31' '6.27'
Set-CellText 'E31' '  +2.54%  '
Set-CellText 'E32' '  +1.59%  '
Set-CellText 'E33' '  -0.60%  '
Set-CellText 'E34' '  +0.00%  '
Set-CellText 'E35' '  +4.89%  '
Set-CellText 'E36' '  -0.45%  '
Set-CellText 'E37' '  +0.12%  '
Set-CellText 'D39' '40.13'
Set-CellText 'E39' '  +2.30%  '
Set-CellText 'E40' '  +12.25%  '
Set-CellText 'D41' '313.66'
Set-CellText 'E41' '  +5.91%  '
Set-CellText 'E42' '  +2.30%  '
Set-CellText 'D43' '142.68'
Set-CellText 'E43' '  -1.53%  '
Set-CellText 'E44' '  +1.63%  '
Set-CellText 'B45' 'Polygon'
Set-CellText 'C45' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText 'D45' '0.422'
Set-CellText 'E45' '  +10.63%  '
Set-CellText 'B46' 'Hedera'
Set-CellText 'C46' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText 'D46' '0.0523'
Set-CellText 'E46' '  +4.12%  '
Set-CellText 'E47' '  -0.03%  '
Set-CellText 'E48' '  +1.91%  '
Set-CellText 'E49' '  +2.93%  '
Set-CellText 'E50' '  -0.27%  '
Set-CellText 'D51' '1.61'
Set-CellText 'E51' '  +4.47%  '
